$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text while writing,
# so numeric-looking values (e.g. "71.90") keep their exact original
# text formatting instead of being coerced into floating point numbers.
$priceCol = $ws.Range("D2:D51")
$priceCol.NumberFormat = "@"

$ws.Range("D2").Value = '42.419.03'
$ws.Range("E2").Value = '  +1.19%  '
$ws.Range("D3").Value = '2.275.16'
$ws.Range("E3").Value = '  +2.50%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '250.78'
$ws.Range("E5").Value = '  -0.02%  '
$ws.Range("D6").Value = '0.633'
$ws.Range("E6").Value = '  +1.54%  '
$ws.Range("D7").Value = '71.90'
$ws.Range("E7").Value = '  +5.83%  '
$ws.Range("E8").Value = '  -0.14%  '
$ws.Range("E9").Value = '  +1.27%  '
$ws.Range("D10").Value = '38.67'
$ws.Range("E10").Value = '  -1.74%  '
$ws.Range("E11").Value = '  +2.65%  '
$ws.Range("D12").Value = '59.03'
$ws.Range("E12").Value = '  -1.38%  '
$ws.Range("D13").Value = '7.32'
$ws.Range("E13").Value = '  +2.90%  '
$ws.Range("D14").Value = '0.105'
$ws.Range("E14").Value = '  +1.79%  '
$ws.Range("D15").Value = '2.615.32'
$ws.Range("E15").Value = '  +2.41%  '
$ws.Range("D16").Value = '14.97'
$ws.Range("E16").Value = '  +2.16%  '
$ws.Range("E17").Value = '  +0.11%  '
$ws.Range("D18").Value = '2.274.98'
$ws.Range("E18").Value = '  +2.56%  '
$ws.Range("D19").Value = '42.349.04'
$ws.Range("E19").Value = '  +1.18%  '
$ws.Range("D20").Value = '0.0₃0995'
$ws.Range("E20").Value = '  +3.25%  '
$ws.Range("D21").Value = '6.28'
$ws.Range("E21").Value = '  +1.05%  '
$ws.Range("D22").Value = '71.85'
$ws.Range("E22").Value = '  -1.10%  '
$ws.Range("D23").Value = '234.83'
$ws.Range("E23").Value = '  +1.19%  '
$ws.Range("D24").Value = '2.21'
$ws.Range("E24").Value = '  +6.86%  '
$ws.Range("E25").Value = '  -0.66%  '
$ws.Range("D26").Value = '11.45'
$ws.Range("E26").Value = '  +0.65%  '
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.34%  '
$ws.Range("D28").Value = '2.43'
$ws.Range("E28").Value = '  +0.36%  '
$ws.Range("E29").Value = '  -1.13%  '
$ws.Range("E30").Value = '  +2.59%  '
$ws.Range("D31").Value = '167.11'
$ws.Range("E31").Value = '  +0.03%  '
$ws.Range("D32").Value = '21.02'
$ws.Range("E32").Value = '  +2.86%  '
$ws.Range("D33").Value = '6.34'
$ws.Range("E33").Value = '  +6.98%  '
$ws.Range("E34").Value = '  +4.38%  '
$ws.Range("D35").Value = '0.0809'
$ws.Range("E35").Value = '  +1.01%  '
$ws.Range("D36").Value = '30.91'
$ws.Range("E36").Value = '  +20.98%  '
$ws.Range("D37").Value = '0.126'
$ws.Range("E37").Value = '  +2.44%  '
$ws.Range("D38").Value = '4.72'
$ws.Range("E38").Value = '  +14.88%  '
$ws.Range("D39").Value = '4.73'
$ws.Range("E39").Value = '  +2.43%  '
$ws.Range("E40").Value = '  -0.37%  '
$ws.Range("D41").Value = '13.80'
$ws.Range("E41").Value = '  +13.07%  '
$ws.Range("D42").Value = '2.33'
$ws.Range("E42").Value = '  +3.78%  '
$ws.Range("D43").Value = '5.89'
$ws.Range("E43").Value = '  +4.49%  '
$ws.Range("D44").Value = '0.213'
$ws.Range("E44").Value = '  +6.54%  '
$ws.Range("D45").Value = '9.15'
$ws.Range("E45").Value = '  +6.44%  '
$ws.Range("D46").Value = '61.25'
$ws.Range("E46").Value = '  -1.48%  '
$ws.Range("E47").Value = '  -3.87%  '
$ws.Range("E48").Value = '  +3.18%  '
$ws.Range("E49").Value = '  +0.08%  '
$ws.Range("E50").Value = '  +0.73%  '
$ws.Range("D51").Value = '97.75'
$ws.Range("E51").Value = '  +4.72%  '

# Restore the original (default/general) cell style for column D so that
# no visible formatting change is introduced by the temporary text format.
$priceCol.Style = "Normal"
